$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# New rows of log data for August 15-16, 2016 (Crestron Logout entries)
$rows = @(
    @{ Row = 225; A = "Crestron Logout"; B = 42597; C = "1700"; D = "TEL"; E = "3069" },
    @{ Row = 226; A = "Crestron Logout"; B = 42597; C = "1700"; D = "TEL"; E = "3072" },
    @{ Row = 227; A = "Crestron Logout"; B = 42597; C = "1630"; D = "TEL"; E = "2116" },
    @{ Row = 231; A = "Crestron Logout"; B = 42598; C = "1700"; D = "TEL"; E = "3069" },
    @{ Row = 232; A = "Crestron Logout"; B = 42598; C = "1700"; D = "TEL"; E = "3072" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}

[void]$ws.Range("E232").Select()
